# Replace the per-row SE value in column D with a categorical "F" label,
# and populate the new sex/environ/sire/dam columns (E:G) with their data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(2, "F", 1, 1, 1),
    @(3, "F", 1, 2, 2),
    @(4, "F", 1, 2, 1),
    @(5, "F", 1, 1, 2),
    @(6, "F", 1, 3, 3),
    @(7, "F", 1, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("E$r").Value = $row[2]
    $ws.Range("F$r").Value = $row[3]
    $ws.Range("G$r").Value = $row[4]
}

# Restore the workbook's saved selection to D2:G7 (with D2 as the active cell).
$ws.Range("D2:G7").Select()
